$d = $word.ActiveDocument

$replacements = @(
    @{Old="71×99=7029"; New="15×51=765"},
    @{Old="63×46=2898"; New="50×95=4750"},
    @{Old="42×21=882";  New="76×17=1292"},
    @{Old="72×51=3672"; New="78×27=2106"},
    @{Old="27×14=378";  New="32×13=416"},
    @{Old="61×58=3538"; New="76×15=1140"},
    @{Old="15×29=435";  New="42×36=1512"},
    @{Old="12×30=360";  New="20×80=1600"},
    @{Old="12×68=816";  New="65×66=4290"},
    @{Old="84×53=4452"; New="43×69=2967"},
    @{Old="94×66=6204"; New="33×91=3003"},
    @{Old="98×95=9310"; New="56×57=3192"},
    @{Old="84×11=924";  New="28×80=2240"},
    @{Old="78×93=7254"; New="49×63=3087"},
    @{Old="95×93=8835"; New="51×20=1020"},
    @{Old="67×84=5628"; New="81×69=5589"},
    @{Old="50×34=1700"; New="82×56=4592"},
    @{Old="85×45=3825"; New="18×60=1080"},
    @{Old="82×19=1558"; New="78×96=7488"},
    @{Old="77×20=1540"; New="64×48=3072"},
    @{Old="11×45=495";  New="81×59=4779"},
    @{Old="66×98=6468"; New="22×60=1320"},
    @{Old="71×40=2840"; New="72×44=3168"},
    @{Old="32×43=1376"; New="84×24=2016"},
    @{Old="20×82=1640"; New="33×94=3102"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
